$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("BH2").Value = "2026-02-28 05:09:57"
$ws.Range("BH3").Value = "2026-02-28 05:09:57"
$ws.Range("BH4").Value = "2026-02-28 05:09:57"
$ws.Range("J5").Value = 3.4
$ws.Range("P5").Value = 1.77
$ws.Range("BH5").Value = "2026-02-28 05:09:57"
$ws.Range("F6").Value = 1.04
$ws.Range("G6").Value = 1000
$ws.Range("H6").Value = 1.04
$ws.Range("I6").Value = 1000
$ws.Range("J6").Value = 1.01
$ws.Range("K6").Value = 1000
$ws.Range("P6").Value = 1.05
$ws.Range("Q6").Value = 1.01
$ws.Range("BH6").Value = "2026-02-28 05:09:57"
$ws.Range("P7").Value = 1.56
$ws.Range("BH7").Value = "2026-02-28 05:09:57"
$ws.Range("BH8").Value = "2026-02-28 05:09:57"
$ws.Range("J9").Value = 3.1
$ws.Range("BH9").Value = "2026-02-28 05:09:57"
$ws.Range("BH10").Value = "2026-02-28 05:09:57"
$ws.Range("J11").Value = 2.54
$ws.Range("BH11").Value = "2026-02-28 05:09:57"
$ws.Range("BH12").Value = "2026-02-28 05:09:57"
$ws.Range("AF13").Value = 23
$ws.Range("AQ13").Value = 8.800000000000001
$ws.Range("AU13").Value = 6.8
$ws.Range("BC13").Value = 34
$ws.Range("BD13").Value = 44
$ws.Range("BF13").Value = 34
$ws.Range("BH13").Value = "2026-02-28 05:09:57"
$ws.Range("J14").Value = 5.8
$ws.Range("O14").Value = 1.33
$ws.Range("Q14").Value = 1.94
$ws.Range("R14").Value = 1.38
$ws.Range("S14").Value = 3.4
$ws.Range("T14").Value = 2.5
$ws.Range("BH14").Value = "2026-02-28 05:09:57"
$ws.Range("H15").Value = 1.62
$ws.Range("I15").Value = 1.67
$ws.Range("Q15").Value = 1.01
$ws.Range("BH15").Value = "2026-02-28 05:09:57"
$ws.Range("BH16").Value = "2026-02-28 05:09:57"
$ws.Range("P17").Value = 1.41
$ws.Range("BH17").Value = "2026-02-28 05:09:57"
$ws.Range("BH18").Value = "2026-02-28 05:09:57"
$ws.Range("BH19").Value = "2026-02-28 05:09:57"
